# Updates Target cluster (column D) assignments and all dependent NATMI-derived
# statistics (columns G-T) for the Rln1-Rxfp1 ligand-receptor pair sheet,
# reflecting a recomputation with updated TPM values ("update scripts wuth new tpm").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value2 = 'FAPs'
$ws.Range('G2').Value2 = [double]"1.212696"
$ws.Range('H2').Value2 = [double]"3.638088"
$ws.Range('I2').Value2 = [double]"0.1802717240630916"
$ws.Range('J2').Value2 = [double]"0.1802717240630916"
$ws.Range('K2').Value2 = [double]"2"
$ws.Range('L2').Value2 = [double]"0.6666666666666666"
$ws.Range('M2').Value2 = [double]"0.029119"
$ws.Range('N2').Value2 = [double]"0.087357"
$ws.Range('O2').Value2 = [double]"0.4059019501247578"
$ws.Range('P2').Value2 = [double]"0.4059019501247578"
$ws.Range('Q2').Value2 = [double]"0.035312494824"
$ws.Range('R2').Value2 = [double]"0.317812453416"
$ws.Range('S2').Value2 = [double]"0.0731726443495611"
$ws.Range('T2').Value2 = [double]"0.07317264434956113"
# Row 3
$ws.Range('D3').Value2 = 'MuSCs'
$ws.Range('G3').Value2 = [double]"1.212696"
$ws.Range('H3').Value2 = [double]"3.638088"
$ws.Range('I3').Value2 = [double]"0.1802717240630916"
$ws.Range('J3').Value2 = [double]"0.1802717240630916"
$ws.Range('M3').Value2 = [double]"0.04229766666666667"
$ws.Range('N3').Value2 = [double]"0.126893"
$ws.Range('O3').Value2 = [double]"0.5896049103927664"
$ws.Range('P3').Value2 = [double]"0.5896049103927664"
$ws.Range('Q3').Value2 = [double]"0.05129421117600001"
$ws.Range('R3').Value2 = [double]"0.4616479005840001"
$ws.Range('S3').Value2 = [double]"0.1062890937125686"
$ws.Range('T3').Value2 = [double]"0.1062890937125687"
# Row 4
$ws.Range('D4').Value2 = 'Neutrophils'
$ws.Range('G4').Value2 = [double]"1.212696"
$ws.Range('H4').Value2 = [double]"3.638088"
$ws.Range('I4').Value2 = [double]"0.1802717240630916"
$ws.Range('J4').Value2 = [double]"0.1802717240630916"
$ws.Range('K4').Value2 = [double]"1"
$ws.Range('L4').Value2 = [double]"0.3333333333333333"
$ws.Range('M4').Value2 = [double]"0.0003223333333333333"
$ws.Range('N4').Value2 = [double]"0.000967"
$ws.Range('O4').Value2 = [double]"0.004493139482475827"
$ws.Range('P4').Value2 = [double]"0.004493139482475827"
$ws.Range('Q4').Value2 = [double]"0.000390892344"
$ws.Range('R4').Value2 = [double]"0.003518031096"
$ws.Range('S4').Value2 = [double]"0.0008099860009618642"
$ws.Range('T4').Value2 = [double]"0.0008099860009618647"
# Row 5
$ws.Range('D5').Value2 = 'FAPs'
$ws.Range('I5').Value2 = [double]"0.3676388369633322"
$ws.Range('J5').Value2 = [double]"0.3676388369633324"
$ws.Range('K5').Value2 = [double]"2"
$ws.Range('L5').Value2 = [double]"0.6666666666666666"
$ws.Range('M5').Value2 = [double]"0.029119"
$ws.Range('N5').Value2 = [double]"0.087357"
$ws.Range('O5').Value2 = [double]"0.4059019501247578"
$ws.Range('P5').Value2 = [double]"0.4059019501247578"
$ws.Range('Q5').Value2 = [double]"0.072014868637"
$ws.Range('R5').Value2 = [double]"0.6481338177330001"
$ws.Range('S5').Value2 = [double]"0.1492253208650144"
$ws.Range('T5').Value2 = [double]"0.1492253208650145"
# Row 6
$ws.Range('D6').Value2 = 'MuSCs'
$ws.Range('I6').Value2 = [double]"0.3676388369633322"
$ws.Range('J6').Value2 = [double]"0.3676388369633324"
$ws.Range('M6').Value2 = [double]"0.04229766666666667"
$ws.Range('N6').Value2 = [double]"0.126893"
$ws.Range('O6').Value2 = [double]"0.5896049103927664"
$ws.Range('P6').Value2 = [double]"0.5896049103927664"
$ws.Range('Q6').Value2 = [double]"0.1046073322796667"
$ws.Range('R6').Value2 = [double]"0.941465990517"
$ws.Range('S6').Value2 = [double]"0.2167616635246664"
$ws.Range('T6').Value2 = [double]"0.2167616635246664"
# Row 7
$ws.Range('D7').Value2 = 'Neutrophils'
$ws.Range('I7').Value2 = [double]"0.3676388369633322"
$ws.Range('J7').Value2 = [double]"0.3676388369633324"
$ws.Range('K7').Value2 = [double]"1"
$ws.Range('L7').Value2 = [double]"0.3333333333333333"
$ws.Range('M7').Value2 = [double]"0.0003223333333333333"
$ws.Range('N7').Value2 = [double]"0.000967"
$ws.Range('O7').Value2 = [double]"0.004493139482475827"
$ws.Range('P7').Value2 = [double]"0.004493139482475827"
$ws.Range('Q7').Value2 = [double]"0.0007971699803333332"
$ws.Range('R7').Value2 = [double]"0.007174529823"
$ws.Range('S7').Value2 = [double]"0.001651852573651441"
$ws.Range('T7').Value2 = [double]"0.001651852573651442"
# Row 8
$ws.Range('D8').Value2 = 'FAPs'
$ws.Range('G8').Value2 = [double]"0.3159016666666667"
$ws.Range('H8').Value2 = [double]"0.947705"
$ws.Range('I8').Value2 = [double]"0.04695994551347087"
$ws.Range('J8').Value2 = [double]"0.04695994551347088"
$ws.Range('K8').Value2 = [double]"2"
$ws.Range('L8').Value2 = [double]"0.6666666666666666"
$ws.Range('M8').Value2 = [double]"0.029119"
$ws.Range('N8').Value2 = [double]"0.087357"
$ws.Range('O8').Value2 = [double]"0.4059019501247578"
$ws.Range('P8').Value2 = [double]"0.4059019501247578"
$ws.Range('Q8').Value2 = [double]"0.009198740631666669"
$ws.Range('R8').Value2 = [double]"0.08278866568500001"
$ws.Range('S8').Value2 = [double]"0.0190611334616702"
$ws.Range('T8').Value2 = [double]"0.0190611334616702"
# Row 9
$ws.Range('D9').Value2 = 'MuSCs'
$ws.Range('G9').Value2 = [double]"0.3159016666666667"
$ws.Range('H9').Value2 = [double]"0.947705"
$ws.Range('I9').Value2 = [double]"0.04695994551347087"
$ws.Range('J9').Value2 = [double]"0.04695994551347088"
$ws.Range('M9').Value2 = [double]"0.04229766666666667"
$ws.Range('N9').Value2 = [double]"0.126893"
$ws.Range('O9').Value2 = [double]"0.5896049103927664"
$ws.Range('P9').Value2 = [double]"0.5896049103927664"
$ws.Range('Q9').Value2 = [double]"0.01336190339611111"
$ws.Range('R9').Value2 = [double]"0.120257130565"
$ws.Range('S9').Value2 = [double]"0.02768781446651918"
$ws.Range('T9').Value2 = [double]"0.02768781446651919"
# Row 10
$ws.Range('D10').Value2 = 'Neutrophils'
$ws.Range('G10').Value2 = [double]"0.3159016666666667"
$ws.Range('H10').Value2 = [double]"0.947705"
$ws.Range('I10').Value2 = [double]"0.04695994551347087"
$ws.Range('J10').Value2 = [double]"0.04695994551347088"
$ws.Range('K10').Value2 = [double]"1"
$ws.Range('L10').Value2 = [double]"0.3333333333333333"
$ws.Range('M10').Value2 = [double]"0.0003223333333333333"
$ws.Range('N10').Value2 = [double]"0.000967"
$ws.Range('O10').Value2 = [double]"0.004493139482475827"
$ws.Range('P10').Value2 = [double]"0.004493139482475827"
$ws.Range('Q10').Value2 = [double]"0.0001018256372222222"
$ws.Range('R10').Value2 = [double]"0.000916430735"
$ws.Range('S10').Value2 = [double]"0.0002109975852814895"
$ws.Range('T10').Value2 = [double]"0.0002109975852814896"
# Row 11
$ws.Range('D11').Value2 = 'FAPs'
$ws.Range('G11').Value2 = [double]"0.1812183333333333"
$ws.Range('H11').Value2 = [double]"0.543655"
$ws.Range('I11').Value2 = [double]"0.02693877227420559"
$ws.Range('J11').Value2 = [double]"0.02693877227420559"
$ws.Range('K11').Value2 = [double]"2"
$ws.Range('L11').Value2 = [double]"0.6666666666666666"
$ws.Range('M11').Value2 = [double]"0.029119"
$ws.Range('N11').Value2 = [double]"0.087357"
$ws.Range('O11').Value2 = [double]"0.4059019501247578"
$ws.Range('P11').Value2 = [double]"0.4059019501247578"
$ws.Range('Q11').Value2 = [double]"0.005276896648333334"
$ws.Range('R11').Value2 = [double]"0.047492069835"
$ws.Range('S11').Value2 = [double]"0.0109345002000668"
$ws.Range('T11').Value2 = [double]"0.01093450020006681"
# Row 12
$ws.Range('D12').Value2 = 'MuSCs'
$ws.Range('G12').Value2 = [double]"0.1812183333333333"
$ws.Range('H12').Value2 = [double]"0.543655"
$ws.Range('I12').Value2 = [double]"0.02693877227420559"
$ws.Range('J12').Value2 = [double]"0.02693877227420559"
$ws.Range('M12').Value2 = [double]"0.04229766666666667"
$ws.Range('N12').Value2 = [double]"0.126893"
$ws.Range('O12').Value2 = [double]"0.5896049103927664"
$ws.Range('P12').Value2 = [double]"0.5896049103927664"
$ws.Range('Q12').Value2 = [double]"0.007665112657222224"
$ws.Range('R12').Value2 = [double]"0.06898601391500001"
$ws.Range('S12').Value2 = [double]"0.01588323241282412"
$ws.Range('T12').Value2 = [double]"0.01588323241282413"
# Row 13
$ws.Range('D13').Value2 = 'Neutrophils'
$ws.Range('G13').Value2 = [double]"0.1812183333333333"
$ws.Range('H13').Value2 = [double]"0.543655"
$ws.Range('I13').Value2 = [double]"0.02693877227420559"
$ws.Range('J13').Value2 = [double]"0.02693877227420559"
$ws.Range('K13').Value2 = [double]"1"
$ws.Range('L13').Value2 = [double]"0.3333333333333333"
$ws.Range('M13').Value2 = [double]"0.0003223333333333333"
$ws.Range('N13').Value2 = [double]"0.000967"
$ws.Range('O13').Value2 = [double]"0.004493139482475827"
$ws.Range('P13').Value2 = [double]"0.004493139482475827"
$ws.Range('Q13').Value2 = [double]"5.841270944444444E-05"
$ws.Range('R13').Value2 = [double]"0.000525714385"
$ws.Range('S13').Value2 = [double]"0.0001210396613146582"
$ws.Range('T13').Value2 = [double]"0.0001210396613146583"
# Row 14
$ws.Range('D14').Value2 = 'FAPs'
$ws.Range('G14').Value2 = [double]"2.500729666666667"
$ws.Range('H14').Value2 = [double]"7.502189"
$ws.Range('I14').Value2 = [double]"0.3717426695773057"
$ws.Range('J14').Value2 = [double]"0.3717426695773058"
$ws.Range('K14').Value2 = [double]"2"
$ws.Range('L14').Value2 = [double]"0.6666666666666666"
$ws.Range('M14').Value2 = [double]"0.029119"
$ws.Range('N14').Value2 = [double]"0.087357"
$ws.Range('O14').Value2 = [double]"0.4059019501247578"
$ws.Range('P14').Value2 = [double]"0.4059019501247578"
$ws.Range('Q14').Value2 = [double]"0.07281874716366668"
$ws.Range('R14').Value2 = [double]"0.6553687244730001"
$ws.Range('S14').Value2 = [double]"0.1508910745260119"
$ws.Range('T14').Value2 = [double]"0.1508910745260119"
# Row 15
$ws.Range('D15').Value2 = 'MuSCs'
$ws.Range('G15').Value2 = [double]"2.500729666666667"
$ws.Range('H15').Value2 = [double]"7.502189"
$ws.Range('I15').Value2 = [double]"0.3717426695773057"
$ws.Range('J15').Value2 = [double]"0.3717426695773058"
$ws.Range('M15').Value2 = [double]"0.04229766666666667"
$ws.Range('N15').Value2 = [double]"0.126893"
$ws.Range('O15').Value2 = [double]"0.5896049103927664"
$ws.Range('P15').Value2 = [double]"0.5896049103927664"
$ws.Range('Q15').Value2 = [double]"0.1057750298641111"
$ws.Range('R15').Value2 = [double]"0.9519752687770001"
$ws.Range('S15').Value2 = [double]"0.2191813033852951"
$ws.Range('T15').Value2 = [double]"0.2191813033852952"
# Row 16
$ws.Range('D16').Value2 = 'Neutrophils'
$ws.Range('G16').Value2 = [double]"2.500729666666667"
$ws.Range('H16').Value2 = [double]"7.502189"
$ws.Range('I16').Value2 = [double]"0.3717426695773057"
$ws.Range('J16').Value2 = [double]"0.3717426695773058"
$ws.Range('K16').Value2 = [double]"1"
$ws.Range('L16').Value2 = [double]"0.3333333333333333"
$ws.Range('M16').Value2 = [double]"0.0003223333333333333"
$ws.Range('N16').Value2 = [double]"0.000967"
$ws.Range('O16').Value2 = [double]"0.004493139482475827"
$ws.Range('P16').Value2 = [double]"0.004493139482475827"
$ws.Range('Q16').Value2 = [double]"0.0008060685292222223"
$ws.Range('R16').Value2 = [double]"0.007254616763"
$ws.Range('S16').Value2 = [double]"0.001670291665998758"
$ws.Range('T16').Value2 = [double]"0.001670291665998758"
# Row 17
$ws.Range('D17').Value2 = 'FAPs'
$ws.Range('G17').Value2 = [double]"0.04337633333333333"
$ws.Range('H17').Value2 = [double]"0.130129"
$ws.Range('I17').Value2 = [double]"0.006448051608593866"
$ws.Range('J17').Value2 = [double]"0.006448051608593868"
$ws.Range('K17').Value2 = [double]"2"
$ws.Range('L17').Value2 = [double]"0.6666666666666666"
$ws.Range('M17').Value2 = [double]"0.029119"
$ws.Range('N17').Value2 = [double]"0.087357"
$ws.Range('O17').Value2 = [double]"0.4059019501247578"
$ws.Range('P17').Value2 = [double]"0.4059019501247578"
$ws.Range('Q17').Value2 = [double]"0.001263075450333333"
$ws.Range('R17').Value2 = [double]"0.011367679053"
$ws.Range('S17').Value2 = [double]"0.002617276722433332"
$ws.Range('T17').Value2 = [double]"0.002617276722433333"
# Row 18
$ws.Range('D18').Value2 = 'MuSCs'
$ws.Range('G18').Value2 = [double]"0.04337633333333333"
$ws.Range('H18').Value2 = [double]"0.130129"
$ws.Range('I18').Value2 = [double]"0.006448051608593866"
$ws.Range('J18').Value2 = [double]"0.006448051608593868"
$ws.Range('M18').Value2 = [double]"0.04229766666666667"
$ws.Range('N18').Value2 = [double]"0.126893"
$ws.Range('O18').Value2 = [double]"0.5896049103927664"
$ws.Range('P18').Value2 = [double]"0.5896049103927664"
$ws.Range('Q18').Value2 = [double]"0.001834717688555555"
$ws.Range('R18').Value2 = [double]"0.016512459197"
$ws.Range('S18').Value2 = [double]"0.00380180289089292"
$ws.Range('T18').Value2 = [double]"0.003801802890892921"
# Row 19
$ws.Range('D19').Value2 = 'Neutrophils'
$ws.Range('G19').Value2 = [double]"0.04337633333333333"
$ws.Range('H19').Value2 = [double]"0.130129"
$ws.Range('I19').Value2 = [double]"0.006448051608593866"
$ws.Range('J19').Value2 = [double]"0.006448051608593868"
$ws.Range('K19').Value2 = [double]"1"
$ws.Range('L19').Value2 = [double]"0.3333333333333333"
$ws.Range('M19').Value2 = [double]"0.0003223333333333333"
$ws.Range('N19').Value2 = [double]"0.000967"
$ws.Range('O19').Value2 = [double]"0.004493139482475827"
$ws.Range('P19').Value2 = [double]"0.004493139482475827"
$ws.Range('Q19').Value2 = [double]"1.398163811111111E-05"
$ws.Range('R19').Value2 = [double]"0.000125834743"
$ws.Range('S19').Value2 = [double]"2.897199526761487E-05"
$ws.Range('T19').Value2 = [double]"2.897199526761488E-05"

Write-Host "Updated Target cluster assignments and NATMI stats for Rln1-Rxfp1 (rows 2-19)"
